$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (match style of existing header cells, e.g. AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in team record data for every data row (2-43)
$lastRow = 43
$wins = $ws.Range("AD2:AD" + $lastRow)
$losses = $ws.Range("AE2:AE" + $lastRow)
$ties = $ws.Range("AF2:AF" + $lastRow)

$wins.Value = 69
$losses.Value = 93
$ties.Value = 0
